# "quick data cleaning; added missing negatives"
#
# Row 18 ("RBC" / "Deposit" / 158.99) was a stray positive deposit entry
# that got removed during cleanup; every row beneath it shifts up by one
# (old row 23 disappears, so the used range shrinks from A1:F23 to A1:F22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Delete()

# Leave the selection where the author left it when they saved.
$ws.Range("A18:XFD18").Select()
